$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet / sheet1.xml) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 66
$ws1.Range("F4").Value = 104
$ws1.Range("F6").Value = 1354
$ws1.Range("F7").Value = 1588
$ws1.Range("F9").Value = 452
$ws1.Range("F13").Value = 75
$ws1.Range("F15").Value = 289
$ws1.Range("F16").Value = 325
$ws1.Range("F25").Value = 4298
$ws1.Range("F28").Value = 1134
$ws1.Range("F31").Value = 656
$ws1.Range("F35").Value = 169
$ws1.Range("F36").Value = 15

# Sheet "全部类型" (4th sheet / sheet4.xml) - column F ("想去人数") updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 66
$ws4.Range("F4").Value = 104
$ws4.Range("F6").Value = 1355
$ws4.Range("F7").Value = 1588
$ws4.Range("F9").Value = 452
$ws4.Range("F13").Value = 75
$ws4.Range("F15").Value = 289
$ws4.Range("F16").Value = 325
$ws4.Range("F25").Value = 4298
$ws4.Range("F28").Value = 1134
$ws4.Range("F31").Value = 656
$ws4.Range("F35").Value = 169
$ws4.Range("F36").Value = 15
